$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Replacement.ClearFormatting()
    $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# "...for tier two, we will use the binary heap. " -> "...for tier two, we will use the binary heap, but can we do better?"
Replace-Text "for tier two, we will use the binary heap. " "for tier two, we will use the binary heap, but can we do better?"

# "...and its lack of Iterators" -> "...and its lack of iterators"
Replace-Text "lack of Iterators" "lack of iterators"

# "...heap it took worst case " -> "...heap took worst case "
Replace-Text "heap it took worst case" "heap took worst case"

# "...its lack for iterator support..." -> "...its lack of iterator support..."
Replace-Text "lack for iterator support" "lack of iterator support"

Write-Output "text edits applied"

# Best-effort: the third chart (Fibonacci VS. priority_queue VS. BinaryHeapWrapper
# at Tier 2) also gained a chart title and axis titles. Attempt this through the
# Word object model; tolerate hosts where chart editing isn't wired up.
try {
    $chartShape = $d.InlineShapes.Item(3)
    $chart = $chartShape.Chart

    $chart.HasTitle = $true
    $chart.ChartTitle.Text = "Fibonacci VS. priority_queue VS. BinaryHeapWrapper at Tier 2"

    $catAxis = $chart.Axes(1)
    $catAxis.HasTitle = $true
    $catAxis.AxisTitle.Text = "Number of Event in heaps"

    $valAxis = $chart.Axes(2)
    $valAxis.HasTitle = $true
    $valAxis.AxisTitle.Text = "Execution time (seconds)"

    Write-Output "chart title edits applied"
} catch {
    Write-Output ("chart title edits skipped: " + $_.Exception.Message)
}
